$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# G4 gets "abs"
$ws.Range("G4").Value = "abs"

# These rows get "x" in column G
$xRows = @(5, 7, 14, 15, 16, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29)
foreach ($r in $xRows) {
    $ws.Range("G$r").Value = "x"
}

# Update the active selection on the sheet to G5
$ws.Activate()
$ws.Range("G5").Select()
